# Jun's Oct 9th updates to MN data
$wb = $excel.ActiveWorkbook

# --- SoCDTtiNTY-psgr: LDVs (row 2) share that is new, B2:H2 -> 0.076 ---
$wsPsgr = $wb.Worksheets.Item("SoCDTtiNTY-psgr")
$wsPsgr.Range("B2:H2").Value = 0.076

# --- SoCDTtiNTY-frgt: HDVs (row 3) share that is new, B3:H3 -> 0.035 ---
$wsFrgt = $wb.Worksheets.Item("SoCDTtiNTY-frgt")
$wsFrgt.Range("B3:H3").Value = 0.035

# Restore view to the first sheet being the active/selected tab
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()
